$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old rows 7-14 (assignment data that no longer applies)
$ws.Rows("7:14").Delete()

# Update the remaining data rows (2-6) with the new values from the
# re-run export (assignments are now read from the folder they live in)
$ws.Range("C2").Value = 111111111
$ws.Range("D2").Value = "test2.txt"
$ws.Range("E2").Value = 70
$ws.Range("F2").Value = 43947.73414990197

$ws.Range("C3").Value = 12345678
$ws.Range("D3").Value = "DocTest.txt"
$ws.Range("E3").Value = 70
$ws.Range("F3").Value = 43947.73425066871

$ws.Range("C4").Value = 161234231
$ws.Range("D4").Value = "test2.txt"
$ws.Range("E4").Value = 30
$ws.Range("F4").Value = 43947.73435902171

$ws.Range("C5").Value = 161234234
$ws.Range("D5").Value = "test2.txt"
$ws.Range("E5").Value = 50
$ws.Range("F5").Value = 43947.73448184071

$ws.Range("C6").Value = 161234236
$ws.Range("D6").Value = "test2.txt"
$ws.Range("E6").Value = 30
$ws.Range("F6").Value = 43947.73462489977
